# "Add files via upload" — the deck was regenerated/re-exported by its
# authoring tool. The only user-visible content change carried by that
# re-export is the "File created on: <date> <time>" stamp on slide 1's
# subtitle placeholder (the rest of the diff is the authoring tool
# minting fresh a16:creationId GUIDs on every shape across all 10
# slides — internal per-shape identifiers that aren't part of the
# PowerPoint object model and aren't settable through automation).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 = title placeholder ("OBESITY PROJECT-..."), unchanged.
# Shape 2 = subtitle placeholder holding the "File created on" stamp.
$stamp = $s.Shapes.Item(2)
$stamp.TextFrame.TextRange.Text = "File created on: 12/10/2023 2:51:46 PM"
